$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $value
    $ws.Range($cell).Style = "Normal"
}

$ws.Range("D2").Value = "27.222.81"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "1.907.43"

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "308.03"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("E6").Value = "  +0.06%  "

Set-TextValue "D7" "0.5251"
$ws.Range("E7").Value = "  +3.36%  "

Set-TextValue "D8" "0.3781"
$ws.Range("E8").Value = "  +3.28%  "

Set-TextValue "D9" "0.07278"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("E10").Value = "  +2.75%  "

Set-TextValue "D11" "0.8978"
$ws.Range("E11").Value = "  +0.83%  "

Set-TextValue "D12" "0.07691"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("D13").Value = "1.905.47"
$ws.Range("E13").Value = "  +1.42%  "

Set-TextValue "D14" "95.18"
$ws.Range("E14").Value = "  +0.46%  "

Set-TextValue "D15" "5.276"
$ws.Range("E15").Value = "  +1.09%  "

Set-TextValue "D17" "0.000008642"
$ws.Range("E17").Value = "  +1.71%  "

Set-TextValue "D18" "14.50"
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "27.280.04"
$ws.Range("E20").Value = "  +1.34%  "

Set-TextValue "D21" "5.077"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").Value = "2.142.02"
$ws.Range("E22").Value = "  +1.27%  "

Set-TextValue "D23" "10.63"
$ws.Range("E23").Value = "  +2.76%  "

Set-TextValue "D24" "6.443"
$ws.Range("E24").Value = "  +1.08%  "

Set-TextValue "D25" "2.330"
$ws.Range("E25").Value = "  +11.72%  "

Set-TextValue "D26" "145.93"
$ws.Range("E26").Value = "  -1.69%  "

Set-TextValue "D27" "1.740"
$ws.Range("E27").Value = "  -1.98%  "

Set-TextValue "D28" "18.13"
$ws.Range("E28").Value = "  +1.62%  "

Set-TextValue "D29" "114.86"
$ws.Range("E29").Value = "  +1.28%  "

Set-TextValue "D30" "4.969"
$ws.Range("E30").Value = "  +5.13%  "

Set-TextValue "D31" "4.820"
$ws.Range("E31").Value = "  +2.77%  "

Set-TextValue "D32" "0.09233"
$ws.Range("E32").Value = "  +1.14%  "

Set-TextValue "D33" "0.8217"
$ws.Range("E33").Value = "  +10.05%  "

Set-TextValue "D34" "0.05072"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("E35").Value = "  +8.00%  "

Set-TextValue "D36" "2.989"
$ws.Range("E36").Value = "  +0.35%  "

Set-TextValue "D37" "3.304"
$ws.Range("E37").Value = "  +2.40%  "

Set-TextValue "D38" "2.601"
$ws.Range("E38").Value = "  +2.82%  "

Set-TextValue "D39" "0.5674"
$ws.Range("E39").Value = "  +1.42%  "

Set-TextValue "D40" "0.01989"
$ws.Range("E40").Value = "  -0.17%  "

Set-TextValue "D41" "1.075"
$ws.Range("E41").Value = "  +0.05%  "

Set-TextValue "D42" "8.995"
$ws.Range("E42").Value = "  +5.22%  "

Set-TextValue "D45" "0.1516"
$ws.Range("E45").Value = "  +2.58%  "

Set-TextValue "D46" "0.4842"
$ws.Range("E46").Value = "  +1.47%  "

Set-TextValue "D47" "10.20"
$ws.Range("E47").Value = "  +1.33%  "

Set-TextValue "D49" "1.628"
$ws.Range("E49").Value = "  +4.66%  "

Set-TextValue "D50" "37.59"
$ws.Range("E50").Value = "  +1.67%  "

Set-TextValue "D51" "63.75"
$ws.Range("E51").Value = "  +1.13%  "

# Row 43 and 44: FraxShare/Quant swap
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "119.37"
$ws.Range("E43").Value = "  +3.30%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "6.632"
$ws.Range("E44").Value = "  +0.27%  "
